$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '73.855.65'
$ws.Range("E2").Value = '  +7.33%  '
$ws.Range("D3").Value = '2.618.85'
$ws.Range("E3").Value = '  +7.47%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '184.26'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +13.72%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '580.47'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +3.71%  '
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("E8").Value = '  +4.14%  '
$ws.Range("E9").Value = '  +17.69%  '
$ws.Range("D10").Value = '2.616.42'
$ws.Range("E10").Value = '  +7.41%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.163'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.24%  '
$ws.Range("E12").Value = '  +8.28%  '
$ws.Range("E13").Value = '  +4.35%  '
$ws.Range("B14").Value = 'WrappedBTC'
$ws.Range("C14").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D14").Value = '73.762.53'
$ws.Range("E14").Value = '  +7.37%  '
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '3.103.09'
$ws.Range("E15").Value = '  +7.51%  '
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000188'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +6.71%  '
$ws.Range("E17").Value = '  +12.02%  '
$ws.Range("D18").Value = '2.618.09'
$ws.Range("E18").Value = '  +7.37%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.14'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +31.90%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.83'
$ws.Range("D20").ClearFormats()
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '371.71'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +9.60%  '
$ws.Range("E22").Value = '  +17.06%  '
$ws.Range("E23").Value = '  +6.05%  '
$ws.Range("E24").Value = '  -0.01%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '69.70'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +4.01%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.13'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +11.22%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.31'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +13.55%  '
$ws.Range("D28").Value = '2.734.74'
$ws.Range("E28").Value = '  +6.55%  '
$ws.Range("E29").Value = '  +0.13%  '
$ws.Range("D30").Value = '0.0₃0934'
$ws.Range("E30").Value = '  +13.81%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '516.78'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +20.69%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.39'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +18.97%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.58'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +6.02%  '
$ws.Range("E34").Value = '  +8.12%  '
$ws.Range("E35").Value = '  +0.05%  '
$ws.Range("E36").Value = '  +12.94%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '161.30'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +1.39%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.16'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +6.50%  '
$ws.Range("E40").Value = '  -0.01%  '
$ws.Range("E41").Value = '  +11.85%  '
$ws.Range("E43").Value = '  +8.49%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '161.62'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +23.99%  '
$ws.Range("E45").Value = '  +9.32%  '
$ws.Range("E46").Value = '  +20.33%  '
$ws.Range("E47").Value = '  +13.86%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '38.59'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +3.19%  '
$ws.Range("E49").Value = '  +8.25%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.527'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +9.62%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '20.46'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +21.25%  '
